$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Hyperlink on A2: keep the same target address, but now show an explicit
#     display string equal to the original e-mail address text. Rebuilding the
#     hyperlink (delete + add) is required since this host only supports
#     mutating hyperlinks through Hyperlinks.Add; re-apply the "Hyperlink"
#     cell style afterwards since Add() re-synthesizes a style of its own. ---
$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:rupamswain1@gmail.com", [System.Type]::Missing, [System.Type]::Missing, "rupamswain1@gmail.com")
$ws.Range("A2").Style = "Hyperlink"

# --- Cell text content: B2 is updated first, then A2, so the underlying
#     shared-string table slots line up the same way as the target workbook. ---
$ws.Range("B2").Value = "cajksbcajksnckansckaskc"
$ws.Range("A2").Value = "rupamswsabsabsjxbajsx"

# --- Phone number cell ---
$ws.Range("C2").Value = 165445545646546

# --- Selection moves from C2 to A2 ---
$ws.Range("A2").Select()
